# Catalog.xlsx update: "AddPriceAgrmnt_RecentOrder" sheet gains four new
# columns (PAItem / NPAItem / ID / UP / Qty data), and the active
# sheet/selection state moves from SmartForm to AddPriceAgrmnt_RecentOrder.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. AddPriceAgrmnt_RecentOrder (sheet3): add columns D:H with new data.
#    New shared-string values are written in the exact order they first
#    appear in the target workbook (NPAI, PAItem, NPAItem, PAI, ID, UP,
#    Qty) so the rebuilt sharedStrings table lines up index-for-index.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("AddPriceAgrmnt_RecentOrder")

$ws3.Range("E2").Value = "NPAI"
$ws3.Range("D1").Value = "PAItem"
$ws3.Range("E1").Value = "NPAItem"
$ws3.Range("D2").Value = "PAI"
$ws3.Range("F1").Value = "ID"
$ws3.Range("G1").Value = "UP"
$ws3.Range("H1").Value = "Qty"
$ws3.Range("F2").Value = "REPOFLOR 100 MG"
$ws3.Range("G2").Value = 2
$ws3.Range("H2").Value = 1

# New column D formatting (matches the width used for the rest of the
# table's custom columns).
$ws3.Columns.Item(4).ColumnWidth = 10.83

# ---------------------------------------------------------------------
# 2. SmartForm (sheet1): selection moves off L6 (and it is no longer the
#    frozen/top-left H1 view) onto C11, and it stops being the active tab.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SmartForm")
$ws1.Range("C11").Select()

# ---------------------------------------------------------------------
# 3. Make AddPriceAgrmnt_RecentOrder the active/tab-selected sheet with
#    its own new selection - done last so it "wins" the active-tab state.
# ---------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("I7").Select()
